$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.473.79"
$ws.Range("E2").Value = "  +1.54%  "

$ws.Range("D3").Value = "3.932.32"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "488.51"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +0.62%  "

$ws.Range("E10").Value = "  +1.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000354"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.07"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.76"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.36%  "

$ws.Range("D14").Value = "4.568.78"
$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.76"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.62%  "

$ws.Range("D16").Value = "3.965.70"
$ws.Range("E16").Value = "  +0.93%  "

$ws.Range("E17").Value = "  -0.61%  "

$ws.Range("E18").Value = "  +0.88%  "

$ws.Range("E19").Value = "  -1.89%  "

$ws.Range("D20").Value = "68.596.23"
$ws.Range("E20").Value = "  +1.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.47"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.93%  "

$ws.Range("E22").Value = "  +4.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.92"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.63%  "

$ws.Range("E24").Value = "  +1.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.39"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +18.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +12.20%  "

$ws.Range("E27").Value = "  +3.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.96"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.78%  "

$ws.Range("E29").Value = "  +1.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "724.85"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.70"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("E32").Value = "  -1.21%  "

$ws.Range("E33").Value = "  +3.05%  "

$ws.Range("D34").Value = "0.0₃0918"
$ws.Range("E34").Value = "  +14.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "42.36"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.18"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +14.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "60.98"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.36%  "

$ws.Range("E38").Value = "  -3.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.399"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +18.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.98"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +14.79%  "

$ws.Range("E42").Value = "  +1.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.17"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.39%  "

$ws.Range("E44").Value = "  +6.11%  "

$ws.Range("E45").Value = "  +0.90%  "

$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("E48").Value = "  -0.81%  "

$ws.Range("E49").Value = "  +1.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "146.02"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("D51").Value = "0.0₆0341"
$ws.Range("E51").Value = "  +36.89%  "
